$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the time_taken column, matching the style of the
# existing header row (B1:E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Per-row timestamps recorded for metadata purposes.
$timestamps = @(
    "2021-10-05 13:39:41.385335",
    "2021-10-05 13:39:41.385348",
    "2021-10-05 13:39:41.385352",
    "2021-10-05 13:39:41.385356",
    "2021-10-05 13:39:41.385359",
    "2021-10-05 13:39:41.385363",
    "2021-10-05 13:39:41.385366"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
